$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new data
$ws.Range("A2").Value = "Take the red block, and move it to the position (0,0.375,0.02). "
$ws.Range("B2").Value = "(0,0.375,0.02)"
$ws.Range("C2").Value = "(-0.25,0.25,0.02)"
$ws.Range("D2").Value = "(0.25,0.5,0.02)"
$ws.Range("E2").Value = "(-0.25,0.5,0.02)"

# Remove old rows 3-5 which are no longer part of the data
$ws.Range("A3:E5").Delete()

# Update selection to match target state
$ws.Range("C10").Select()
